$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.930.60'
$ws.Cells.Item(2, 5).Value = '  -5.76%  '
$ws.Cells.Item(3, 4).Value = '3.037.75'
$ws.Cells.Item(3, 5).Value = '  -6.84%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '554.40'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(5, 5).Value = '  -6.48%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '140.20'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(6, 5).Value = '  -9.34%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 4).Value = '3.035.50'
$ws.Cells.Item(8, 5).Value = '  -6.73%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.482'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(9, 5).Value = '  -12.12%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.155'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(10, 5).Value = '  -11.81%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '6.10'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(11, 5).Value = '  -10.96%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.456'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(12, 5).Value = '  -10.34%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '34.91'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(13, 5).Value = '  -10.24%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000219'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(14, 5).Value = '  -11.36%  '
$ws.Cells.Item(15, 4).Value = '3.523.50'
$ws.Cells.Item(15, 5).Value = '  -6.47%  '
$ws.Cells.Item(16, 4).Value = '63.921.27'
$ws.Cells.Item(16, 5).Value = '  -5.56%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.110'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(17, 5).Value = '  -3.87%  '
$ws.Cells.Item(18, 4).Value = '3.032.60'
$ws.Cells.Item(18, 5).Value = '  -6.52%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '6.51'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(19, 5).Value = '  -10.53%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '477.87'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(20, 5).Value = '  -12.47%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '13.38'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(21, 5).Value = '  -12.57%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.666'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(22, 5).Value = '  -13.45%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '7.10'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(23, 5).Value = '  -9.78%  '
$ws.Cells.Item(24, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '12.32'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(24, 5).Value = '  -10.01%  '
$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '77.09'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(25, 5).Value = '  -10.63%  '
$ws.Cells.Item(26, 5).Value = '  -0.06%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(27, 5).Value = '  -16.02%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '2.06'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(28, 5).Value = '  -3.74%  '
$ws.Cells.Item(29, 2).Value = 'RenderToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '7.53'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(29, 5).Value = '  -8.58%  '
$ws.Cells.Item(30, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(30, 5).Value = '  +0.05%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '25.84'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(31, 5).Value = '  -13.07%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '2.59'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(32, 5).Value = '  -5.02%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '1.10'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(33, 5).Value = '  -4.86%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '493.12'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(34, 5).Value = '  -10.73%  '
$ws.Cells.Item(35, 2).Value = 'NEARProtocol'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '5.25'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(35, 5).Value = '  -9.98%  '
$ws.Cells.Item(36, 2).Value = 'OKB'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '51.87'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(36, 5).Value = '  -4.05%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '5.78'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(37, 5).Value = '  -13.29%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.0402'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(38, 5).Value = '  -10.50%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.0778'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(39, 5).Value = '  -9.08%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(40, 5).Value = '  -8.70%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '8.16'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(41, 5).Value = '  -12.16%  '
$ws.Cells.Item(42, 4).Value = '2.862.69'
$ws.Cells.Item(42, 5).Value = '  -3.10%  '
$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '2.45'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(43, 5).Value = '  -8.08%  '
$ws.Cells.Item(44, 2).Value = 'USDe'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(44, 5).Value = '  -0.17%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.239'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(45, 5).Value = '  -9.42%  '
$ws.Cells.Item(46, 4).Value = '0.0₃0527'
$ws.Cells.Item(46, 5).Value = '  -11.81%  '
$ws.Cells.Item(47, 2).Value = 'Fetch.AI'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(47, 5).Value = '  -7.13%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '24.17'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(48, 5).Value = '  -8.45%  '
$ws.Cells.Item(49, 2).Value = 'Monero'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '116.37'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(49, 5).Value = '  -7.45%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.106'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(50, 5).Value = '  -7.19%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '2.01'
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Cells.Item(51, 5).Value = '  -16.68%  '
